$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap "Periodo Mora" (column E) and "Valor Mora" (column F) between rows 16 and 17
$e16 = $ws.Range("E16").Value()
$e17 = $ws.Range("E17").Value()
$f16 = $ws.Range("F16").Value()
$f17 = $ws.Range("F17").Value()

$ws.Range("E16").Value = $e17
$ws.Range("E17").Value = $e16
$ws.Range("F16").Value = $f17
$ws.Range("F17").Value = $f16
